$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column U ("Storage scan Qty" / "&=result.sort_qty") entirely.
# Columns V and W shift left to U and V.
$ws.Range("U1").EntireColumn.Delete()

# Update the active selection to match the post-edit state.
$ws.Range("P13").Select()
